$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold text-formatted values (prices with "." as thousands
# separator, and percentage strings) in the source data. Force the Text number
# format first so Excel keeps our assignments as literal strings instead of
# auto-converting them to numbers (which would lose exact formatting like
# trailing zeros, e.g. "1.00" -> 1).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.866.61'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '3.056.10'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '554.71'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').Value = '141.65'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.056.16'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  +4.76%  '
$ws.Range('D11').Value = '5.97'
$ws.Range('E11').Value = '  -8.37%  '
$ws.Range('D12').Value = '0.471'
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('D13').Value = '0.0000230'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '34.62'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '3.558.46'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '64.024.18'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.110'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.050.76'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '6.69'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '477.48'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').Value = '13.95'
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('D22').Value = '0.672'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').Value = '7.49'
$ws.Range('E23').Value = '  +3.72%  '
$ws.Range('D24').Value = '14.11'
$ws.Range('E24').Value = '  +11.28%  '
$ws.Range('D25').Value = '81.12'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '2.78'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '7.92'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '2.04'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('D30').Value = '0.995'
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').Value = '26.08'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').Value = '1.14'
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('D33').Value = '2.44'
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').Value = '5.60'
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('D35').Value = '6.15'
$ws.Range('E35').Value = '  +3.62%  '
$ws.Range('D36').Value = '54.87'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').Value = '0.0406'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '2.93'
$ws.Range('E38').Value = '  +16.59%  '
$ws.Range('D39').Value = '440.47'
$ws.Range('E39').Value = '  -4.80%  '
$ws.Range('D40').Value = '0.0804'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').Value = '2.960.65'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').Value = '8.17'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').Value = '0.113'
$ws.Range('E43').Value = '  -4.59%  '
$ws.Range('D44').Value = '28.19'
$ws.Range('E44').Value = '  +3.11%  '
$ws.Range('D45').Value = '0.257'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = '2.12'
$ws.Range('E47').Value = '  +4.69%  '
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '117.47'
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0514'
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('E51').Value = '  +0.57%  '
